{"js": "const pairs = [\n  [\"42-25=\", \"89-12=\"],\n  [\"87-52=\", \"53-37=\"],\n  [\"26+71=\", \"33+4=\"],\n  [\"40+24=\", \"17+18=\"],\n  [\"48-6=\", \"73-27=\"],\n  [\"28+27=\", \"33+16=\"],\n  [\"24+40=\", \"39+28=\"],\n  [\"74-58=\", \"5-3=\"],\n  [\"51-32=\", \"97-36=\"],\n  [\"97-32=\", \"95-91=\"],\n  [\"21+62=\", \"9+29=\"],\n  [\"24+43=\", \"88+1=\"],\n  [\"4+67=\", \"68+9=\"],\n  [\"34+24=\", \"5+81=\"],\n  [\"41-37=\", \"35-25=\"],\n  [\"40-8=\", \"72-29=\"],\n  [\"7+82=\", \"97-25=\"],\n  [\"49-34=\", \"97-0=\"],\n  [\"46-38=\", \"62+19=\"],\n  [\"21+7=\", \"54+17=\"],\n  [\"28+30=\", \"47+3=\"],\n  [\"45+28=\", \"55-17=\"],\n  [\"5+49=\", \"69+30=\"],\n  [\"70-27=\", \"29+51=\"],\n  [\"39+23=\", \"4+53=\"],\n  [\"35+60=\", \"77-51=\"],\n  [\"42+48=\", \"20+54=\"],\n  [\"68-61=\", \"5+18=\"],\n  [\"25+16=\", \"26-21=\"],\n  [\"10+25=\", \"89-46=\"],\n  [\"11+72=\", \"22+29=\"],\n  [\"45-13=\", \"68+20=\"],\n  [\"89-21=\", \"67-51=\"],\n  [\"9+8=\", \"66+12=\"],\n  [\"77-16=\", \"36+38=\"],\n  [\"38+24=\", \"6+56=\"],\n  [\"60-32=\", \"20+46=\"],\n  [\"97-21=\", \"87-20=\"],\n  [\"5+13=\", \"96-78=\"],\n  [\"13+27=\", \"51-43=\"],\n  [\"3+81=\", \"71-57=\"],\n  [\"31+27=\", \"58-54=\"],\n  [\"89-70=\", \"82-81=\"],\n  [\"30+3=\", \"25+32=\"],\n  [\"60-12=\", \"13+40=\"],\n  [\"10-7=\", \"77-13=\"],\n  [\"74-72=\", \"74+4=\"],\n  [\"35+28=\", \"13+26=\"],\n  [\"67+15=\", \"33+45=\"],\n  [\"85-78=\", \"85-34=\"],\n  [\"69-67=\", \"23-16=\"],\n  [\"98-46=\", \"30+54=\"],\n  [\"76-56=\", \"65-25=\"],\n  [\"13+72=\", \"20-10=\"],\n  [\"52+18=\", \"33+62=\"],\n  [\"91-6=\", \"27-8=\"],\n  [\"74-16=\", \"50+19=\"],\n  [\"53+0=\", \"4+2=\"],\n  [\"17-8=\", \"87-84=\"],\n  [\"54-36=\", \"54-48=\"],\n  [\"72-24=\", \"21+25=\"],\n  [\"12+47=\", \"89-36=\"],\n  [\"57-17=\", \"74-11=\"],\n  [\"38+11=\", \"80-42=\"],\n  [\"60-21=\", \"92-21=\"],\n  [\"17+3=\", \"63+9=\"],\n  [\"32+38=\", \"64-41=\"],\n  [\"93-45=\", \"40+33=\"],\n  [\"45-31=\", \"74+13=\"],\n  [\"39-4=\", \"31+55=\"],\n  [\"42-4=\", \"61-20=\"],\n  [\"33+58=\", \"69+6=\"],\n  [\"61-10=\", \"14+3=\"],\n  [\"81-36=\", \"29-26=\"],\n  [\"92-5=\", \"63-38=\"],\n  [\"1+75=\", \"92-19=\"],\n  [\"27+19=\", \"72-0=\"],\n  [\"45+46=\", \"0+67=\"],\n  [\"84-37=\", \"25-6=\"],\n  [\"62+28=\", \"98-11=\"],\n  [\"23+25=\", \"5+18=\"],\n  [\"20+23=\", \"86-3=\"],\n  [\"67+1=\", \"64+3=\"],\n  [\"0+57=\", \"39-31=\"],\n  [\"24+46=\", \"55+22=\"],\n  [\"50-47=\", \"95-21=\"],\n  [\"25+39=\", \"23+51=\"],\n  [\"71-53=\", \"49+21=\"],\n  [\"15+47=\", \"59-31=\"],\n  [\"0+14=\", \"27-26=\"],\n  [\"87-17=\", \"35+7=\"],\n  [\"42+0=\", \"81+8=\"],\n  [\"3+37=\", \"89-38=\"],\n  [\"95-82=\", \"68-19=\"],\n  [\"24+15=\", \"74+5=\"],\n  [\"61-3=\", \"56+8=\"],\n  [\"34+35=\", \"66-56=\"],\n  [\"18+34=\", \"47+36=\"],\n  [\"23+32=\", \"91-76=\"],\n  [\"22-8=\", \"40+26=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"42-25=\", \"89-12=\"),\n    @(\"87-52=\", \"53-37=\"),\n    @(\"26+71=\", \"33+4=\"),\n    @(\"40+24=\", \"17+18=\"),\n    @(\"48-6=\", \"73-27=\"),\n    @(\"28+27=\", \"33+16=\"),\n    @(\"24+40=\", \"39+28=\"),\n    @(\"74-58=\", \"5-3=\"),\n    @(\"51-32=\", \"97-36=\"),\n    @(\"97-32=\", \"95-91=\"),\n    @(\"21+62=\", \"9+29=\"),\n    @(\"24+43=\", \"88+1=\"),\n    @(\"4+67=\", \"68+9=\"),\n    @(\"34+24=\", \"5+81=\"),\n    @(\"41-37=\", \"35-25=\"),\n    @(\"40-8=\", \"72-29=\"),\n    @(\"7+82=\", \"97-25=\"),\n    @(\"49-34=\", \"97-0=\"),\n    @(\"46-38=\", \"62+19=\"),\n    @(\"21+7=\", \"54+17=\"),\n    @(\"28+30=\", \"47+3=\"),\n    @(\"45+28=\", \"55-17=\"),\n    @(\"5+49=\", \"69+30=\"),\n    @(\"70-27=\", \"29+51=\"),\n    @(\"39+23=\", \"4+53=\"),\n    @(\"35+60=\", \"77-51=\"),\n    @(\"42+48=\", \"20+54=\"),\n    @(\"68-61=\", \"5+18=\"),\n    @(\"25+16=\", \"26-21=\"),\n    @(\"10+25=\", \"89-46=\"),\n    @(\"11+72=\", \"22+29=\"),\n    @(\"45-13=\", \"68+20=\"),\n    @(\"89-21=\", \"67-51=\"),\n    @(\"9+8=\", \"66+12=\"),\n    @(\"77-16=\", \"36+38=\"),\n    @(\"38+24=\", \"6+56=\"),\n    @(\"60-32=\", \"20+46=\"),\n    @(\"97-21=\", \"87-20=\"),\n    @(\"5+13=\", \"96-78=\"),\n    @(\"13+27=\", \"51-43=\"),\n    @(\"3+81=\", \"71-57=\"),\n    @(\"31+27=\", \"58-54=\"),\n    @(\"89-70=\", \"82-81=\"),\n    @(\"30+3=\", \"25+32=\"),\n    @(\"60-12=\", \"13+40=\"),\n    @(\"10-7=\", \"77-13=\"),\n    @(\"74-72=\", \"74+4=\"),\n    @(\"35+28=\", \"13+26=\"),\n    @(\"67+15=\", \"33+45=\"),\n    @(\"85-78=\", \"85-34=\"),\n    @(\"69-67=\", \"23-16=\"),\n    @(\"98-46=\", \"30+54=\"),\n    @(\"76-56=\", \"65-25=\"),\n    @(\"13+72=\", \"20-10=\"),\n    @(\"52+18=\", \"33+62=\"),\n    @(\"91-6=\", \"27-8=\"),\n    @(\"74-16=\", \"50+19=\"),\n    @(\"53+0=\", \"4+2=\"),\n    @(\"17-8=\", \"87-84=\"),\n    @(\"54-36=\", \"54-48=\"),\n    @(\"72-24=\", \"21+25=\"),\n    @(\"12+47=\", \"89-36=\"),\n    @(\"57-17=\", \"74-11=\"),\n    @(\"38+11=\", \"80-42=\"),\n    @(\"60-21=\", \"92-21=\"),\n    @(\"17+3=\", \"63+9=\"),\n    @(\"32+38=\", \"64-41=\"),\n    @(\"93-45=\", \"40+33=\"),\n    @(\"45-31=\", \"74+13=\"),\n    @(\"39-4=\", \"31+55=\"),\n    @(\"42-4=\", \"61-20=\"),\n    @(\"33+58=\", \"69+6=\"),\n    @(\"61-10=\", \"14+3=\"),\n    @(\"81-36=\", \"29-26=\"),\n    @(\"92-5=\", \"63-38=\"),\n    @(\"1+75=\", \"92-19=\"),\n    @(\"27+19=\", \"72-0=\"),\n    @(\"45+46=\", \"0+67=\"),\n    @(\"84-37=\", \"25-6=\"),\n    @(\"62+28=\", \"98-11=\"),\n    @(\"23+25=\", \"5+18=\"),\n    @(\"20+23=\", \"86-3=\"),\n    @(\"67+1=\", \"64+3=\"),\n    @(\"0+57=\", \"39-31=\"),\n    @(\"24+46=\", \"55+22=\"),\n    @(\"50-47=\", \"95-21=\"),\n    @(\"25+39=\", \"23+51=\"),\n    @(\"71-53=\", \"49+21=\"),\n    @(\"15+47=\", \"59-31=\"),\n    @(\"0+14=\", \"27-26=\"),\n    @(\"87-17=\", \"35+7=\"),\n    @(\"42+0=\", \"81+8=\"),\n    @(\"3+37=\", \"89-38=\"),\n    @(\"95-82=\", \"68-19=\"),\n    @(\"24+15=\", \"74+5=\"),\n    @(\"61-3=\", \"56+8=\"),\n    @(\"34+35=\", \"66-56=\"),\n    @(\"18+34=\", \"47+36=\"),\n    @(\"23+32=\", \"91-76=\"),\n    @(\"22-8=\", \"40+26=\"),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}"}
